# Atualizado por script em 07-01-2024 08:45
#
# The source scraper re-ran and produced a slightly different ordering /
# extra rows for the Morocco Botola Pro 2023-2024 sheet:
#   - Match pairs that used to appear in one order now appear swapped
#     (rows 12/13, 15/16, 47/48 keep their Indice/date but the match
#     details in columns F:V trade places).
#   - Two newly scraped matches are appended as rows 112/113.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($Sheet, [int]$Row, [int]$StartCol, [object[]]$Values)
    for ($i = 0; $i -lt $Values.Count; $i++) {
        $Sheet.Cells.Item($Row, $StartCol + $i).Value = $Values[$i]
    }
}

function Swap-MatchDetails {
    # Swap the F:V ("match details") payload between two rows, leaving
    # A (Indice) and B:E (pais/torneio/temporada/data_partida) untouched.
    param($Sheet, [int]$RowA, [int]$RowB)

    $colCount = 17 # F..V inclusive
    $valsA = @()
    $valsB = @()

    for ($i = 0; $i -lt $colCount; $i++) {
        $valsA += $Sheet.Cells.Item($RowA, 6 + $i).Value2
        $valsB += $Sheet.Cells.Item($RowB, 6 + $i).Value2
    }

    Set-RowValues $Sheet $RowA 6 $valsB
    Set-RowValues $Sheet $RowB 6 $valsA
}

# --- Swap the three re-ordered match pairs -------------------------------
Swap-MatchDetails $ws 12 13
Swap-MatchDetails $ws 15 16
Swap-MatchDetails $ws 47 48

# --- Append the two newly scraped matches --------------------------------

# Copy the number-formatting / border styling used by the existing data
# rows (bold bordered "Indice" cell in column A, datetime-formatted cell
# in column E) onto the two new rows before filling in their values.
$ws.Range("A111").Copy() | Out-Null
$ws.Range("A112").PasteSpecial(-4122) | Out-Null
$ws.Range("A113").PasteSpecial(-4122) | Out-Null

$ws.Range("E111").Copy() | Out-Null
$ws.Range("E112").PasteSpecial(-4122) | Out-Null
$ws.Range("E113").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

Set-RowValues $ws 112 1 @(
    111,
    "morocco",
    "botola-pro",
    "2023-2024",
    45297.75,
    "Raja Casablanca",
    2,
    "Mouloudia Oujda",
    0,
    1.19,
    "04/01/2024 23:42",
    1.17,
    "06/01/2024 17:58",
    5.41,
    "04/01/2024 23:42",
    6.41,
    "06/01/2024 17:58",
    13.34,
    "04/01/2024 23:42",
    17.68,
    "06/01/2024 17:58",
    "https://www.betexplorer.com/football/morocco/botola-pro/raja-casablanca-mouloudia-oujda/roV8VqeC/"
)

Set-RowValues $ws 113 1 @(
    112,
    "morocco",
    "botola-pro",
    "2023-2024",
    45297.83333333334,
    "Chabab Mohammedia",
    1,
    "IR Tanger",
    0,
    2.44,
    "04/01/2024 23:42",
    2.63,
    "06/01/2024 19:50",
    2.79,
    "04/01/2024 23:42",
    3.15,
    "06/01/2024 19:52",
    2.94,
    "04/01/2024 23:42",
    2.68,
    "06/01/2024 19:52",
    "https://www.betexplorer.com/football/morocco/botola-pro/chabab-mohammedia-ir-tanger/4fW4WPu6/"
)
